# Update "想去人数" (number of interested attendees) figures for several
# events across the workbook's sheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F8").Value = 6982
$ws.Range("F14").Value = 275
$ws.Range("F18").Value = 4465
$ws.Range("F29").Value = 8038
$ws.Range("F31").Value = 1387
$ws.Range("F38").Value = 1627
$ws.Range("F40").Value = 939
$ws.Range("F42").Value = 4082

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 29

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 238

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 238
$ws.Range("F12").Value = 6982
$ws.Range("F18").Value = 275
$ws.Range("F21").Value = 4465
$ws.Range("F29").Value = 8038
$ws.Range("F31").Value = 1387
$ws.Range("F38").Value = 1627
$ws.Range("F40").Value = 939
$ws.Range("F42").Value = 4082
